# Applies bold "label:" prefixes to bullet paragraphs on slides 3, 4 and 9,
# and inserts a new "username," column name into the Features bullet on
# slide 4 (CSV feature list).
#
# Strategy: every affected paragraph keeps its existing concatenated text
# (so we never touch a:p boundaries) - we just locate the label substring
# with IndexOf and flip Font.Bold = $true on that TextRange.Characters(...)
# slice, which causes PowerPoint to split the run in two (bold label +
# normal remainder), matching the target OOXML.
#
# NOTE: this COM-interop host mis-parses a parenthesized expression used
# directly as a function-call argument (e.g. `Foo $x ("bar".Length)`
# throws "Cannot invoke a value ... not a COM collection"). Always assign
# such expressions to a variable first, then pass the variable.

$p = $ppt.ActivePresentation

function Set-BoldRange($rangeObj, $startPos1, $len) {
    $r = $rangeObj.Characters($startPos1, $len)
    $r.Font.Bold = $true
}

function Set-BoldLabel($textRange, $full, $label) {
    $idx0 = $full.IndexOf($label)
    $start1 = $idx0 + 1
    $len = $label.Length
    Set-BoldRange $textRange $start1 $len
}

# ---------------------------------------------------------------------
# Slide 3 - Business Understanding bullets
# ---------------------------------------------------------------------
$s3 = $p.Slides.Item(3)
$tr3 = $s3.Shapes.Item(2).TextFrame.TextRange
$full3 = $tr3.Text
Set-BoldLabel $tr3 $full3 "Customer:"
Set-BoldLabel $tr3 $full3 "Problem:"
Set-BoldLabel $tr3 $full3 "Why it matters:"
Set-BoldLabel $tr3 $full3 "Goal: "

# ---------------------------------------------------------------------
# Slide 4 - Data Understanding bullets
# ---------------------------------------------------------------------
$s4 = $p.Slides.Item(4)
$tr4 = $s4.Shapes.Item(2).TextFrame.TextRange

# Insert the new "username," feature name right before the existing
# "captions, hashtags, likes, followers" text.
$full4 = $tr4.Text
$oldSuffix = "captions, hashtags, likes, followers"
$newSuffix = "username, captions, hashtags, likes, followers"
$capIdx0 = $full4.IndexOf($oldSuffix)
$capStart1 = $capIdx0 + 1
$capLen = $oldSuffix.Length
$capRange = $tr4.Characters($capStart1, $capLen)
$capRange.Text = $newSuffix

# Re-read text after the edit before computing any further offsets.
$full4 = $tr4.Text
Set-BoldLabel $tr4 $full4 "Dataset:"
Set-BoldLabel $tr4 $full4 "Features: "

# Bold the single space that separates "username," from "captions,...".
$userLabel = "username,"
$userIdx0 = $full4.IndexOf($userLabel)
$spaceIdx0 = $userIdx0 + $userLabel.Length
$spaceStart1 = $spaceIdx0 + 1
Set-BoldRange $tr4 $spaceStart1 1

Set-BoldLabel $tr4 $full4 "Issues:"

# ---------------------------------------------------------------------
# Slide 9 - Modeling bullets
# ---------------------------------------------------------------------
$s9 = $p.Slides.Item(9)
$tr9 = $s9.Shapes.Item(2).TextFrame.TextRange
$full9 = $tr9.Text
Set-BoldLabel $tr9 $full9 "Target: "
Set-BoldLabel $tr9 $full9 "Baseline:"
Set-BoldLabel $tr9 $full9 "Improved:"
Set-BoldLabel $tr9 $full9 "Features: "
